# Add a new "Wind_Farm" object mapping row to the Object_Mapping sheet,
# mirroring the commit's integration of a windfarm object into the mapping.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Object_Mapping")

$ws.Range("A14").Value = "Wind_Farm"
$ws.Range("B14").Value = "Wind_farm"

# Leave selection where Excel would land after typing into the last row.
$ws.Range("A15").Select()
